$d = $word.ActiveDocument

# Fix "relevent" -> "relevant"
$d.Content.Find.Execute("relevent", $true, $false, $false, $false, $false,
                         $true, 1, $false, "relevant", 2)

# Fix "digestable" -> "digestible"
$d.Content.Find.Execute("digestable", $true, $false, $false, $false, $false,
                         $true, 1, $false, "digestible", 2)
